$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $dateRange = $ws.Range("A2:A8")

    # Force text interpretation so the "yyyy/mm/dd" strings aren't
    # auto-converted into date serials, then strip the number format
    # back off so the cells keep their original (General/default) style.
    $dateRange.NumberFormat = "@"

    $ws.Range("A2").Value = "2025/12/31"
    $ws.Range("A3").Value = "2024/12/31"
    $ws.Range("A4").Value = "2023/12/31"
    $ws.Range("A5").Value = "2022/12/31"
    $ws.Range("A6").Value = "2021/12/31"
    $ws.Range("A7").Value = "2020/12/31"
    $ws.Range("A8").Value = "2015/12/31"

    $dateRange.ClearFormats()

    $ws.Range("B8").Value = "Upgrade"
    $ws.Range("C8").Value = "Upgrade"
    $ws.Range("D8").Value = "Upgrade"
    $ws.Range("E8").Value = "Upgrade"
}
